$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 238, pushing existing rows 238-284 down to 239-285
$ws.Rows.Item(238).Insert()

# Fill in the new row 238 with the inserted record's data
$ws.Cells.Item(238, 1).Value = 8
$ws.Cells.Item(238, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(238, 3).Value = "Coquimbo"
$ws.Cells.Item(238, 4).Value = 44637
$ws.Cells.Item(238, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(238, 5).Value = 4
$ws.Cells.Item(238, 6).Value = 100112032
$ws.Cells.Item(238, 7).Value = "Zapallo italiano"
$ws.Cells.Item(238, 8).Value = "Sin especificar"
$ws.Cells.Item(238, 9).Value = "Primera"
$ws.Cells.Item(238, 10).Value = 560
$ws.Cells.Item(238, 11).Value = 11000
$ws.Cells.Item(238, 12).Value = 12000
$ws.Cells.Item(238, 13).Value = 11500
$ws.Cells.Item(238, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(238, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(238, 16).Value = 192
$ws.Cells.Item(238, 17).Value = 60
$ws.Cells.Item(238, 18).Value = "Hortaliza"
